$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its values stored as text, matching the
# source data which uses dotted thousand separators (e.g. "22.410.77") and
# values like "1.000" that Excel would otherwise coerce to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (D) column updates ---
$ws.Range("D2").Value = "22.410.77"
$ws.Range("D3").Value = "1.562.98"
$ws.Range("D6").Value = "285.80"
$ws.Range("D7").Value = "0.3637"
$ws.Range("D8").Value = "48.30"
$ws.Range("D9").Value = "0.3340"
$ws.Range("D11").Value = "0.07414"
$ws.Range("D12").Value = "1.000"
$ws.Range("D13").Value = "20.80"
$ws.Range("D14").Value = "5.929"
$ws.Range("D15").Value = "6.886"
$ws.Range("D16").Value = "1.563.81"
$ws.Range("D18").Value = "88.53"
$ws.Range("D19").Value = "0.06697"
$ws.Range("D20").Value = "1.000"
$ws.Range("D21").Value = "6.337"
$ws.Range("D22").Value = "16.08"
$ws.Range("D24").Value = "22.405.15"
$ws.Range("D25").Value = "2.415"
$ws.Range("D26").Value = "2.549"
$ws.Range("D27").Value = "149.70"
$ws.Range("D28").Value = "19.38"
$ws.Range("D29").Value = "4.999"
$ws.Range("D30").Value = "123.00"
$ws.Range("D31").Value = "1.737.31"
$ws.Range("D32").Value = "1.065"
$ws.Range("D33").Value = "6.120"
$ws.Range("D34").Value = "1.997"
$ws.Range("D35").Value = "9.601"
$ws.Range("D36").Value = "0.08222"
$ws.Range("D37").Value = "0.02393"
$ws.Range("D38").Value = "1.302"
$ws.Range("D39").Value = "0.06388"
$ws.Range("D40").Value = "0.2205"
$ws.Range("D41").Value = "5.338"
$ws.Range("D42").Value = "11.14"
$ws.Range("D43").Value = "0.6076"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D45").Value = "13.73"
$ws.Range("D47").Value = "0.5745"
$ws.Range("D48").Value = "2.010"
$ws.Range("D49").Value = "124.86"
$ws.Range("D50").Value = "1.211"
$ws.Range("D51").Value = "0.07214"

# --- Coin name (B), Link (C) and Volume(1h) (E) column updates ---
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  -2.74%  "
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("E18").Value = "  -2.66%  "
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E38").Value = "  -5.45%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  -1.57%  "
